$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "During the update:" -> "After the update:"
# ------------------------------------------------------------------
$d.Content.Find.Execute("During the update:", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "After the update:", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Delete the now-unwanted paragraphs (in reverse document order so
#    earlier indices/ranges stay valid while we work).
#    Targets (by their exact current text):
#      - "Separate color for rows..." -> handled separately (text replace)
#      - "Think about the MPCs..."                         (delete)
#      - "Do we want deflators and MPCs..."                (delete)
#      - "Create a monthly Haver pull sheet..." (2nd/standalone copy)  (delete)
#      - "After the update:" (the duplicate heading paragraph)        (delete)
#      - "Taxes Louise "                                   (delete)
#      - "Medicare (why are BEA and CBO so different in 2020)"        (delete)
#      - "State purchases "                                (delete)
# ------------------------------------------------------------------

function Remove-ParagraphByText($doc, [string]$exactText) {
    for ($i = $doc.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        # Strip the trailing paragraph mark character for comparison.
        if ($t.Length -gt 0) {
            $lastChar = $t.Substring($t.Length - 1)
            if ($lastChar -eq [char]13 -or $lastChar -eq [char]7) {
                $t = $t.Substring(0, $t.Length - 1)
            }
        }
        if ($t -eq $exactText) {
            $p.Range.Delete() | Out-Null
            return $true
        }
    }
    return $false
}

Remove-ParagraphByText $d "Think about the MPCs – whether to edit in spreadsheet or in the code, how the functions should work" | Out-Null
Remove-ParagraphByText $d "Do we want deflators and MPCs in the spreadsheet? Seems risky for getting out of date given that we update them in the code and would have to update in two places. " | Out-Null
Remove-ParagraphByText $d "Create a monthly Haver pull sheet to call from (after the update)" | Out-Null
Remove-ParagraphByText $d "After the update:" | Out-Null
Remove-ParagraphByText $d "Taxes Louise " | Out-Null
Remove-ParagraphByText $d "Medicare (why are BEA and CBO so different in 2020)" | Out-Null
Remove-ParagraphByText $d "State purchases " | Out-Null

# ------------------------------------------------------------------
# 3. Rewrite the (formerly) strike-through paragraph: its text is
#    replaced entirely, and the new text is split across three runs
#    with different strike formatting.
# ------------------------------------------------------------------

function Find-ParagraphIndexByText($doc, [string]$exactText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Length -gt 0) {
            $lastChar = $t.Substring($t.Length - 1)
            if ($lastChar -eq [char]13 -or $lastChar -eq [char]7) {
                $t = $t.Substring(0, $t.Length - 1)
            }
        }
        if ($t -eq $exactText) {
            return $i
        }
    }
    return -1
}

$idx = Find-ParagraphIndexByText $d "When update comes in, remember to go over these re cbo"
$p2 = $d.Paragraphs.Item($idx)
$pStart = $p2.Range.Start
$pEnd = $p2.Range.End

# Remove the old text (but not the paragraph mark itself, so the
# paragraph's own pPr/rPr -- the strike default -- stays put).
$textRange = $d.Range($pStart, $pEnd - 1)
$textRange.Delete() | Out-Null

$part1 = "Create a monthly Haver pull sheet to call from (after the update)"
$part2 = " "
$part3 = "(I think this actually more complicated than just copy pasting it in)"

$insertPoint = $d.Range($pStart, $pStart)
$insertPoint.InsertBefore($part1 + $part2 + $part3) | Out-Null

$run1 = $d.Range($pStart, $pStart + $part1.Length)
$run1.Font.StrikeThrough = 1

$run2 = $d.Range($pStart + $part1.Length, $pStart + $part1.Length + $part2.Length)
$run2.Font.StrikeThrough = 1

$run3Start = $pStart + $part1.Length + $part2.Length
$run3End = $run3Start + $part3.Length
$run3 = $d.Range($run3Start, $run3End)
$run3.Font.StrikeThrough = 0

# ------------------------------------------------------------------
# 4. "Separate color for rows..." -> "Master script for Louise"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Separate color for rows that shouldn’t be dragged/need Louise’s input", `
                         $false, $false, $false, $false, $false, $true, 1, $false, `
                         "Master script for Louise", 2) | Out-Null
